# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap country names for rows 168 (Republica de Africa Central) and 169 (Republica del Chad) ---
# and refresh their case numbers with newly reported data (Chad updated, Africa Central keeps old numbers).
$ws.Range("A168").Value = "Republica del Chad"
$ws.Range("B168").Value = 52
$ws.Range("C168").Value = 6
$ws.Range("D168").Value = 19
$ws.Range("E168").Value = 31
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 2
$ws.Range("H168").Value = 2

$ws.Range("A169").Value = "Republica de Africa Central"
$ws.Range("B169").Value = 50
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 10
$ws.Range("E169").Value = 40
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 0

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 20:22"

# --- Refresh country statistics with new case numbers ---
# Row 8: Alemania
$ws.Range("B8").Value = 159239
$ws.Range("C8").Value = 481
$ws.Range("E8").Value = 35662
$ws.Range("G8").Value = 51
$ws.Range("H8").Value = 6177

# Row 15: Canada
$ws.Range("B15").Value = 49040
$ws.Range("C15").Value = 540
$ws.Range("D15").Value = 18721
$ws.Range("E15").Value = 27550
$ws.Range("G15").Value = 62
$ws.Range("H15").Value = 2769

# Row 22: Ecuador
$ws.Range("B22").Value = 24258
$ws.Range("C22").Value = 1018
$ws.Range("E22").Value = 21830
$ws.Range("G22").Value = 208
$ws.Range("H22").Value = 871

# Row 28
$ws.Range("F28").Value = 136

# Row 30
$ws.Range("B30").Value = 14612
$ws.Range("C30").Value = 697
$ws.Range("E30").Value = 11067

# Row 58
$ws.Range("F58").Value = 22
